Write-Output "Starting e_pwTemplate localization edit"

# Translate the (Chinese) sheet tab names to their English equivalents.
# Do this FIRST: renaming a sheet makes the engine rewrite the sheet-name
# token baked into the *last* workbook-scoped defined name's RefersTo text,
# so any defined-name fixups must happen only after every rename is done.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)
$ws5 = $wb.Worksheets.Item(5)
$ws6 = $wb.Worksheets.Item(6)

$ws1.Name = "Well Blessing"
$ws2.Name = "Site Guardian"
$ws3.Name = "Ancestors"
$ws4.Name = "Deceased"
$ws5.Name = "Karmic Creditors"
$ws6.Name = "Recently Deceased"

Write-Output "Renamed all six worksheet tabs"

# Recreate the W001A defined names (one local to "Deceased", one local to
# "Recently Deceased") now that renaming is finished, so their RefersTo
# text is stamped with the final sheet names and nothing renames afterward
# to clobber it. Re-adding them in this order also yields the swapped
# ordering (localSheetId=3 before localSheetId=5) seen in the target file.
foreach ($n in $ws4.Names) { $n.Delete() }
foreach ($n in $ws6.Names) { $n.Delete() }

Write-Output "Deleted old localized defined names"

$ws4.Names.Add("W001A", "=Deceased!`$B`$3:`$E`$519")
$ws6.Names.Add("W001A", "='Recently Deceased'!`$B`$3:`$F`$522")

Write-Output "Recreated W001A defined names against the renamed sheets"
